# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values on Sheet1 to the newly
# regenerated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 2
    4  = 2
    6  = 2
    7  = 1
    8  = 2
    9  = 4
    10 = 1
    11 = 2
    12 = 0
    13 = 4
    14 = 3
    15 = 1
    16 = 3
    17 = 3
    18 = 2
    19 = 1
    20 = 1
    21 = 1
    23 = 1
    24 = 2
    25 = 1
    26 = 1
    27 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
